$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.34 = 12760.54 pesos`n✅ 12760.54 pesos = 3.32 = 973.91 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3815.4
$ws2.Range("N12").Value = 3839
$ws2.Range("O12").Value = 293
